# Update the team-specific transition-probability matrix on Sheet1 with
# refreshed counts (more games simulated -> recalculated frequencies).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2085106382978723
$ws.Range("C2").Value = 0.5531914893617021
$ws.Range("J2").Value = 0.01702127659574468
$ws.Range("P2").Value = 0.1276595744680851
$ws.Range("S2").Value = 0.09361702127659574

# Row 3
$ws.Range("B3").Value = 0.007633587786259542
$ws.Range("C3").Value = 0.02290076335877863
$ws.Range("J3").Value = 0.03053435114503817
$ws.Range("P3").Value = 0.7251908396946565
$ws.Range("S3").Value = 0.2137404580152672

# Row 4
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("P4").Value = 0.6458333333333334
$ws.Range("S4").Value = 0.2708333333333333

# Row 6
$ws.Range("B6").Value = 0.03626943005181347
$ws.Range("D6").Value = 0.02590673575129534
$ws.Range("F6").Value = 0.04663212435233161
$ws.Range("J6").Value = 0.1865284974093264
$ws.Range("O6").Value = 0.03626943005181347
$ws.Range("Q6").Value = 0.155440414507772
$ws.Range("R6").Value = 0.06735751295336788
$ws.Range("S6").Value = 0.4455958549222798

# Row 7
$ws.Range("B7").Value = 0.1173184357541899
$ws.Range("D7").Value = 0.02793296089385475
$ws.Range("F7").Value = 0.05586592178770949
$ws.Range("J7").Value = 0.1452513966480447
$ws.Range("O7").Value = 0.0223463687150838
$ws.Range("Q7").Value = 0.1843575418994413
$ws.Range("R7").Value = 0.05027932960893855
$ws.Range("S7").Value = 0.3966480446927375

# Row 8
$ws.Range("B8").Value = 0.07526881720430108
$ws.Range("D8").Value = 0.01720430107526882
$ws.Range("F8").Value = 0.05806451612903226
$ws.Range("J8").Value = 0.1182795698924731
$ws.Range("O8").Value = 0.02580645161290323
$ws.Range("Q8").Value = 0.1913978494623656
$ws.Range("R8").Value = 0.0989247311827957
$ws.Range("S8").Value = 0.4150537634408602

# Row 9
$ws.Range("B9").Value = 0.09881422924901186
$ws.Range("D9").Value = 0.007905138339920948
$ws.Range("F9").Value = 0.04347826086956522
$ws.Range("J9").Value = 0.1739130434782609
$ws.Range("O9").Value = 0.007905138339920948
$ws.Range("Q9").Value = 0.1778656126482213
$ws.Range("R9").Value = 0.08300395256916997
$ws.Range("S9").Value = 0.4071146245059288

# Row 10
$ws.Range("B10").Value = 0.1064400715563506
$ws.Range("D10").Value = 0.0259391771019678
$ws.Range("E10").Value = 0.0008944543828264759
$ws.Range("F10").Value = 0.07871198568872988
$ws.Range("J10").Value = 0.08944543828264759
$ws.Range("O10").Value = 0.01610017889087657
$ws.Range("Q10").Value = 0.2164579606440072
$ws.Range("R10").Value = 0.08407871198568873
$ws.Range("S10").Value = 0.3819320214669052

# Row 11
$ws.Range("G11").Value = 0.1485507246376812
$ws.Range("J11").Value = 0.1050724637681159
$ws.Range("K11").Value = 0.1956521739130435
$ws.Range("L11").Value = 0.5326086956521739
$ws.Range("S11").Value = 0.01811594202898551

# Row 12
$ws.Range("G12").Value = 0.7383720930232558
$ws.Range("J12").Value = 0.1569767441860465
$ws.Range("L12").Value = 0.02325581395348837
$ws.Range("S12").Value = 0.08139534883720931

# Row 13
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.08333333333333333

# Row 15
$ws.Range("F15").Value = 0.004950495049504951
$ws.Range("H15").Value = 0.2128712871287129
$ws.Range("I15").Value = 0.09405940594059406
$ws.Range("J15").Value = 0.301980198019802
$ws.Range("K15").Value = 0.0594059405940594
$ws.Range("M15").Value = 0.0198019801980198
$ws.Range("O15").Value = 0.02475247524752475
$ws.Range("S15").Value = 0.2821782178217822

# Row 16
$ws.Range("F16").Value = 0.006535947712418301
$ws.Range("H16").Value = 0.2026143790849673
$ws.Range("I16").Value = 0.0718954248366013
$ws.Range("J16").Value = 0.3594771241830065
$ws.Range("K16").Value = 0.1176470588235294
$ws.Range("M16").Value = 0.006535947712418301
$ws.Range("O16").Value = 0.03267973856209151
$ws.Range("S16").Value = 0.2026143790849673

# Row 17
$ws.Range("F17").Value = 0.006880733944954129
$ws.Range("H17").Value = 0.213302752293578
$ws.Range("I17").Value = 0.1055045871559633
$ws.Range("J17").Value = 0.3922018348623853
$ws.Range("K17").Value = 0.07568807339449542
$ws.Range("M17").Value = 0.02064220183486239
$ws.Range("O17").Value = 0.04587155963302753
$ws.Range("S17").Value = 0.1399082568807339

# Row 18
$ws.Range("H18").Value = 0.1758241758241758
$ws.Range("I18").Value = 0.1043956043956044
$ws.Range("J18").Value = 0.3736263736263736
$ws.Range("K18").Value = 0.1318681318681319
$ws.Range("M18").Value = 0.01648351648351648
$ws.Range("N18").Value = 0.005494505494505495
$ws.Range("O18").Value = 0.06593406593406594
$ws.Range("S18").Value = 0.1263736263736264

# Row 19
$ws.Range("F19").Value = 0.01722846441947565
$ws.Range("H19").Value = 0.1932584269662921
$ws.Range("I19").Value = 0.1161048689138577
$ws.Range("J19").Value = 0.3363295880149813
$ws.Range("K19").Value = 0.09288389513108614
$ws.Range("M19").Value = 0.0149812734082397
$ws.Range("N19").Value = 0.002247191011235955
$ws.Range("O19").Value = 0.06591760299625468
$ws.Range("S19").Value = 0.1610486891385768
